# Saldo_guide.xlsx update
# - Rename the sheet/workbook tab to "Saldo_guide"
# - Shift every "Dt. Referencia" (col G) from 2024-07-26 (45499) to 2024-07-29 (45502)
# - Recalculate "Vl. Projetado" (col D) / "Saldo Previsto" (col E) / "Vl. Total" (col H)
#   for the rows whose balances moved
# - Leave the active selection on D13, matching the author's last cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new workbook/file name
$ws.Name = "Saldo_guide"

# Every data row's reference date moved from 45499 (2024-07-26) to 45502 (2024-07-29)
$lastRow = 274
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45502
}

# Update "Vl. Projetado" (D) values that changed from 0 to a new figure
$ws.Cells.Item(5, 4).Value = 6913.57
$ws.Cells.Item(8, 4).Value = 2882.18
$ws.Cells.Item(15, 4).Value = 5396.91
$ws.Cells.Item(17, 4).Value = 4010.95
$ws.Cells.Item(43, 4).Value = 6073.07
$ws.Cells.Item(49, 4).Value = 1938.06
$ws.Cells.Item(58, 4).Value = 390.46
$ws.Cells.Item(60, 4).Value = 6481.06
$ws.Cells.Item(99, 4).Value = 5342.66
$ws.Cells.Item(104, 4).Value = 18657.77
$ws.Cells.Item(108, 4).Value = 16833.22
$ws.Cells.Item(132, 4).Value = 2283.97
$ws.Cells.Item(143, 4).Value = 18939.65
$ws.Cells.Item(158, 4).Value = 524.93
$ws.Cells.Item(172, 4).Value = 979.36
$ws.Cells.Item(173, 4).Value = 12456.58
$ws.Cells.Item(235, 4).Value = 4740.38
$ws.Cells.Item(249, 4).Value = 4848.8
$ws.Cells.Item(264, 4).Value = 13743.97
$ws.Cells.Item(265, 4).Value = 8381.63
$ws.Cells.Item(270, 4).Value = 5205.63
$ws.Cells.Item(271, 4).Value = 7173.16
$ws.Cells.Item(273, 4).Value = 4710.8

# Update "Saldo Previsto" (E) values that were recalculated
$ws.Cells.Item(107, 5).Value = 894.98
$ws.Cells.Item(112, 5).Value = 21.64
$ws.Cells.Item(113, 5).Value = 0

# Update "Vl. Total" (H) = Vl. Projetado + Saldo Previsto for affected rows
$ws.Cells.Item(5, 8).Value = 7505.95
$ws.Cells.Item(8, 8).Value = 3297.74
$ws.Cells.Item(15, 8).Value = 6214.89
$ws.Cells.Item(17, 8).Value = 4427.76
$ws.Cells.Item(43, 8).Value = 7692.23
$ws.Cells.Item(49, 8).Value = 2134.2
$ws.Cells.Item(58, 8).Value = 590.91
$ws.Cells.Item(60, 8).Value = 6818.84
$ws.Cells.Item(99, 8).Value = 5703.15
$ws.Cells.Item(104, 8).Value = 19655.01
$ws.Cells.Item(107, 8).Value = 894.98
$ws.Cells.Item(108, 8).Value = 18028.86
$ws.Cells.Item(112, 8).Value = 21.64
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(132, 8).Value = 2501.04
$ws.Cells.Item(143, 8).Value = 20298.62
$ws.Cells.Item(158, 8).Value = 1222.83
$ws.Cells.Item(172, 8).Value = 1846.51
$ws.Cells.Item(173, 8).Value = 13517.34
$ws.Cells.Item(235, 8).Value = 5199
$ws.Cells.Item(249, 8).Value = 5100.81
$ws.Cells.Item(264, 8).Value = 15639.47
$ws.Cells.Item(265, 8).Value = 9294.44
$ws.Cells.Item(270, 8).Value = 5760.35
$ws.Cells.Item(271, 8).Value = 7952.12
$ws.Cells.Item(273, 8).Value = 5499.66

# Match the author's final cursor position/selection
$ws.Range("D13").Select()
